# edit.ps1 -- applies the commit's two logical changes:
#   1. The table on slide 5 gets a new table style GUID.
#   2. The presentation's live theme (colour scheme) is swapped from the
#      "Integral" / Red-Violet palette to the stock Office palette (the
#      deck's two theme parts effectively trade their colour schemes).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style id change on slide 5's table (B1 - Types of financial
#    documents), shape 2.
# ---------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{2969C1AA-31BE-4BB2-BC22-FF4294999901}")
}

# ---------------------------------------------------------------------
# 2) Theme colour swap: the deck's active theme (currently the
#    "Integral" / Red Violet scheme) takes on the plain "Office Theme"
#    colours.
# ---------------------------------------------------------------------
function HexToBgr([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme colour scheme, in the standard 12-slot theme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToBgr $officeColors[$i - 1]
}
